# Insert 3 new data rows at the top of the Chirimoya price block (row 176),
# pushing the existing rows 176..268 down to 179..271, and populate the
# 3 new rows with the new weekly price entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(176).Resize(3).Insert()

# Row 176
$ws.Cells.Item(176,1).Value = 6
$ws.Cells.Item(176,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(176,3).Value = "Metropolitana"
$ws.Cells.Item(176,4).Value = 44830
$ws.Cells.Item(176,5).Value = 13
$ws.Cells.Item(176,6).Value = "Fruta"
$ws.Cells.Item(176,7).Value = 100107
$ws.Cells.Item(176,8).Value = "Otros"
$ws.Cells.Item(176,9).Value = 100107002
$ws.Cells.Item(176,10).Value = "Chirimoya"
$ws.Cells.Item(176,11).Value = "Cultivar IV Región"
$ws.Cells.Item(176,12).Value = "Especial"
$ws.Cells.Item(176,13).Value = 250
$ws.Cells.Item(176,14).Value = 25000
$ws.Cells.Item(176,15).Value = 25000
$ws.Cells.Item(176,16).Value = 25000
$ws.Cells.Item(176,17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(176,18).Value = "Provincia de Limarí"
$ws.Cells.Item(176,19).Value = 3125
$ws.Cells.Item(176,20).Value = 8

# Row 177
$ws.Cells.Item(177,1).Value = 6
$ws.Cells.Item(177,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(177,3).Value = "Metropolitana"
$ws.Cells.Item(177,4).Value = 44830
$ws.Cells.Item(177,5).Value = 13
$ws.Cells.Item(177,6).Value = "Fruta"
$ws.Cells.Item(177,7).Value = 100107
$ws.Cells.Item(177,8).Value = "Otros"
$ws.Cells.Item(177,9).Value = 100107002
$ws.Cells.Item(177,10).Value = "Chirimoya"
$ws.Cells.Item(177,11).Value = "Cultivar IV Región"
$ws.Cells.Item(177,12).Value = "Primera"
$ws.Cells.Item(177,13).Value = 250
$ws.Cells.Item(177,14).Value = 22000
$ws.Cells.Item(177,15).Value = 22000
$ws.Cells.Item(177,16).Value = 22000
$ws.Cells.Item(177,17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(177,18).Value = "Provincia de Limarí"
$ws.Cells.Item(177,19).Value = 2750
$ws.Cells.Item(177,20).Value = 8

# Row 178
$ws.Cells.Item(178,1).Value = 6
$ws.Cells.Item(178,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(178,3).Value = "Metropolitana"
$ws.Cells.Item(178,4).Value = 44830
$ws.Cells.Item(178,5).Value = 13
$ws.Cells.Item(178,6).Value = "Fruta"
$ws.Cells.Item(178,7).Value = 100107
$ws.Cells.Item(178,8).Value = "Otros"
$ws.Cells.Item(178,9).Value = 100107002
$ws.Cells.Item(178,10).Value = "Chirimoya"
$ws.Cells.Item(178,11).Value = "Cultivar IV Región"
$ws.Cells.Item(178,12).Value = "Segunda"
$ws.Cells.Item(178,13).Value = 250
$ws.Cells.Item(178,14).Value = 18000
$ws.Cells.Item(178,15).Value = 18000
$ws.Cells.Item(178,16).Value = 18000
$ws.Cells.Item(178,17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(178,18).Value = "Provincia de Limarí"
$ws.Cells.Item(178,19).Value = 2250
$ws.Cells.Item(178,20).Value = 8
